# Applies the changes described in the diff for
# ConfigurationFiles/2-Score_weight_PL.xlsx:
#   - new "Weight"/"Group" values for the Ceramide FA1/FA2 scoring rows
#   - the sheet's last active selection moved to E20
#   - the workbook window width setting changed

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (FA1_[FA-H]-): Weight 35 -> 47.5
$ws.Range("B2").Value = 47.5

# Row 3 (FA2_[FA-H]-): Weight 35 -> 47.5
$ws.Range("B3").Value = 47.5

# Row 4 (FA1_[LPL-H]-): Weight 9 -> 1.5
$ws.Range("B4").Value = 1.5

# Row 5 (FA2_[LPL-H]-): Weight 9 -> 1.5
$ws.Range("B5").Value = 1.5

# Row 6 (FA1_[LPL-H2O-H]-): Weight 6 -> 1, Group 2 -> 3
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 3

# Row 7 (FA2_[LPL-H2O-H]-): Weight 6 -> 1, Group 2 -> 3
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 3

# Move the active/selected cell shown when the sheet is reopened
$ws.Range("E20").Select()

# Resize the workbook window (windowWidth 22260 -> 12960)
$excel.ActiveWindow.Width = 12960
